$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T11:42:05.701+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T11:50:01.626+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T11:59:19.102+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T12:21:42.961+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T12:22:54.587+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T12:26:44.353+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
    ,@('testPostVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testGetVideoById', 'Passed', 'GET https://www.videogamedb.uk:443/api/videogame/1', '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}')
    ,@('testPutVideo', 'Passed', '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}', '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}')
    ,@('testDeleteVideo', 'Passed', 'DELETE https://www.videogamedb.uk:443/api/videogame/1', 'Video game deleted')
    ,@('testInvalidPostVideo', 'Passed', '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}', '{"timestamp":"2025-01-02T12:28:15.332+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}')
)

$startRow = 32
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}

Write-Output "Added $($rows.Count) rows"
